$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value  = 2419.2191551329238
$ws.Range("C7").Value  = 0.12777098831577546
$ws.Range("D7").Value  = 3.7370505550693451

$ws.Range("B8").Value  = 2812.9596952487791
$ws.Range("C8").Value  = 0.13282417709981051

$ws.Range("B11").Value = 1078.4121933398442
$ws.Range("D11").Value = 186.80031295621546

$ws.Range("B12").Value = 1132.8490994067702
$ws.Range("D12").Value = 139.0570551314251

$ws.Range("B14").Value = 1091.2984859118972
$ws.Range("C14").Value = 0.061751677021090731
$ws.Range("D14").Value = 115.19004111822326

$ws.Range("B17").Value = 5632.1745357733398
$ws.Range("C17").Value = 0.29248934600350512
$ws.Range("D17").Value = 6.9941655313570523

$ws.Range("B18").Value = 5500.3486655212346
$ws.Range("C18").Value = 0.29603958705060746
$ws.Range("D18").Value = 141.12735216743064

$ws.Range("B19").Value = 1452.694086178293
$ws.Range("C19").Value = 0.13768899938936063
$ws.Range("D19").Value = 146.4778276970068

$ws.Range("B20").Value = 4028.700214152826
$ws.Range("C20").Value = 0.22911358696760745
$ws.Range("D20").Value = 8.7137901629899872

$ws.Range("B22").Value = 1163.0651754465221
$ws.Range("C22").Value = 0.078527036547241039
$ws.Range("D22").Value = 123.93220346095688

$ws.Range("B23").Value = 1066.9938966080795
$ws.Range("C23").Value = 0.19130997388736096
$ws.Range("D23").Value = 159.38187662210208

$ws.Range("B24").Value = 1143.3239908483383
$ws.Range("D24").Value = 193.9725682632789

$ws.Range("B25").Value = 1125.6210000564597
$ws.Range("C25").Value = 0.1012485157420745
$ws.Range("D25").Value = 175.16591353332882

$ws.Range("B26").Value = 2558.367947046987
$ws.Range("C26").Value = 0.17057046589393166
$ws.Range("D26").Value = 13.259241766734227

$ws.Range("B27").Value = 5219.2500087002445
$ws.Range("C27").Value = 0.28416905940597675
$ws.Range("D27").Value = 11.492125457871889

$ws.Range("B28").Value = 5217.2576971377221
$ws.Range("C28").Value = 0.28176715132997271
$ws.Range("D28").Value = 10.277121048361494
